$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2 through 416). Bump it from 45181 (2023-09-12) to 45182 (2023-09-13)
# for every row, matching the diff which updates each row's C value by +1.
for ($row = 2; $row -le 416; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45181) {
        $cell.Value2 = 45182
    }
}
